$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" values replacing the old "Strike#" values in column G (rows 2-27)
$kValues = @{
    2  = 3
    3  = 5
    4  = 4
    5  = 5
    6  = 8
    7  = 3
    8  = 0
    9  = 6
    10 = 4
    11 = 4
    12 = 10
    13 = 5
    14 = 6
    15 = 6
    16 = 5
    17 = 8
    18 = 12
    19 = 6
    20 = 8
    21 = 4
    22 = 4
    23 = 2
    24 = 7
    25 = 3
    26 = 1
    27 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
